$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge row-1 merged ranges so individual cells can be written
$ws.Range("A1:C1").UnMerge()
$ws.Range("D1:R1").UnMerge()
$ws.Range("S1:BK1").UnMerge()

# Row 1
$ws.Range("S1").Value = "Product Features"
$ws.Range("T1").Value = ""
$ws.Range("BL1").Value = "Survey Data"
$ws.Range("BN1").Style = "Normal"

# Row 2
$ws.Range("S2").Value = "Speed Importance Level"
$ws.Range("T2").Value = "Speed Satisfaction Level"
$ws.Range("U2").Value = "Speed Fulfillment Capacity"
$ws.Range("V2").Value = "Charging time Importance Level"
$ws.Range("W2").Value = "Charging time Satisfaction Level"
$ws.Range("X2").Value = "Charging time Fulfillment Capacity"
$ws.Range("Y2").Value = "Eco-Friendly Importance Level"
$ws.Range("Z2").Value = "Eco-Friendly Satisfaction Level"
$ws.Range("AA2").Value = "Eco-Friendly Fulfillment Capacity"
$ws.Range("AB2").Value = "Mileage Importance Level"
$ws.Range("AC2").Value = "Mileage Satisfaction Level"
$ws.Range("AD2").Value = "Mileage Fulfillment Capacity"
$ws.Range("AE2").Value = "Seat capacity Importance Level"
$ws.Range("AF2").Value = "Seat capacity Satisfaction Level"
$ws.Range("AG2").Value = "Seat capacity Fulfillment Capacity"
$ws.Range("AH2").Value = "Brand Importance Level"
$ws.Range("AI2").Value = "Brand Satisfaction Level"
$ws.Range("AJ2").Value = "Brand Fulfillment Capacity"
$ws.Range("AK2").Value = "Design Importance Level"
$ws.Range("AL2").Value = "Design Satisfaction Level"
$ws.Range("AM2").Value = "Design Fulfillment Capacity"
$ws.Range("AN2").Value = "Battery Importance Level"
$ws.Range("AO2").Value = "Battery Satisfaction Level"
$ws.Range("AP2").Value = "Battery Fulfillment Capacity"
$ws.Range("AQ2").Value = "Weight Importance Level"
$ws.Range("AR2").Value = "Weight Satisfaction Level"
$ws.Range("AS2").Value = "Weight Fulfillment Capacity"
$ws.Range("AT2").Value = "Safety Importance Level"
$ws.Range("AU2").Value = "Safety Satisfaction Level"
$ws.Range("AV2").Value = "Safety Fulfillment Capacity"
$ws.Range("AW2").Value = "Price Importance Level"
$ws.Range("AX2").Value = "Price Satisfaction Level"
$ws.Range("AY2").Value = "Price Fulfillment Capacity"
$ws.Range("AZ2").Value = "Maintenance Importance Level"
$ws.Range("BA2").Value = "Maintenance Satisfaction Level"
$ws.Range("BB2").Value = "Maintenance Fulfillment Capacity"
$ws.Range("BC2").Value = "Social value Importance Level"
$ws.Range("BD2").Value = "Social value Satisfaction Level"
$ws.Range("BE2").Value = "Social value Fulfillment Capacity"
$ws.Range("BF2").Value = "Re-sell value Importance Level"
$ws.Range("BG2").Value = "Re-sell value Satisfaction Level"
$ws.Range("BH2").Value = "Re-sell value Fulfillment Capacity"
$ws.Range("BI2").Value = "Overall satisfaction Importance Level"
$ws.Range("BJ2").Value = "Overall satisfaction Satisfaction Level"
$ws.Range("BK2").Value = "Overall satisfaction Fulfillment Capacity"
$ws.Range("BL2").Value = "Age"
$ws.Range("BM2").Value = "Buy Vehicle in Future"

# Row 3
$ws.Range("C3").Value = "sports"
$ws.Range("D3").Value = "'4"
$ws.Range("E3").Value = "'9"
$ws.Range("F3").Value = "'3"
$ws.Range("G3").Value = "'8"
$ws.Range("H3").Value = "'10"
$ws.Range("I3").Value = "'14"
$ws.Range("J3").Value = "'15"
$ws.Range("K3").Value = "'11"
$ws.Range("L3").Value = "'2"
$ws.Range("M3").Value = "'1"
$ws.Range("N3").Value = "'12"
$ws.Range("O3").Value = "'7"
$ws.Range("P3").Value = "'6"
$ws.Range("Q3").Value = "'5"
$ws.Range("R3").Value = "'13"
$ws.Range("S3").Value = "'2"
$ws.Range("T3").Value = "'6"
$ws.Range("U3").Value = "'7"
$ws.Range("Y3").Value = "'5"
$ws.Range("Z3").Value = "'6"
$ws.Range("AA3").Value = "'4"
$ws.Range("AB3").Value = "'3"
$ws.Range("AC3").Value = "'7"
$ws.Range("AD3").Value = "'5"
$ws.Range("AE3").Value = "'2"
$ws.Range("AF3").Value = "'3"
$ws.Range("AI3").Value = "'7"
$ws.Range("AJ3").Value = "'5"
$ws.Range("AK3").Value = "'6"
$ws.Range("AL3").Value = "'6"
$ws.Range("AM3").Value = "'2"
$ws.Range("AN3").Value = "'6"
$ws.Range("AO3").Value = "'7"
$ws.Range("AP3").Value = "'2"
$ws.Range("AU3").Value = "'3"
$ws.Range("AV3").Value = "'1"
$ws.Range("AW3").Value = "'4"
$ws.Range("AX3").Value = "'5"
$ws.Range("AY3").Value = "'6"
$ws.Range("AZ3").Value = "'6"
$ws.Range("BA3").Value = "'7"
$ws.Range("BC3").Value = "'4"
$ws.Range("BE3").Value = "'2"
$ws.Range("BF3").Value = "'3"
$ws.Range("BG3").Value = "'3"
$ws.Range("BI3").Value = "'5"
$ws.Range("BJ3").Value = "'5"
$ws.Range("BK3").Value = "'5"
$ws.Range("BL3").Value = "26-35"
$ws.Range("BM3").Value = 0

# Re-merge row-1 ranges (including the new BL1:BM1 merge)
$ws.Range("A1:C1").Merge()
$ws.Range("D1:R1").Merge()
$ws.Range("S1:BK1").Merge()
$ws.Range("BL1:BM1").Merge()
